# Auto-generated edit script applying the cryptos.xlsx price/volume refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    # Force text interpretation so numeric-looking strings (e.g. "16.30")
    # keep their exact formatting/trailing zeros instead of being coerced
    # into a Double by the smart-entry parser.
    $rng.NumberFormat = "@"
    $rng.Value = $text
    # Drop the explicit text-format style again so the cell keeps the same
    # (default) style index it had before - only its stored text changes.
    $rng.ClearFormats()
}

Set-TextValue "D2" "66.845.66"
Set-TextValue "D3" "3.112.98"
Set-TextValue "E3" "  +0.59%  "
Set-TextValue "E4" "  +0.00%  "
Set-TextValue "D5" "577.21"
Set-TextValue "E5" "  -0.52%  "
Set-TextValue "D6" "171.02"
Set-TextValue "E6" "  +1.87%  "
Set-TextValue "E7" "  +0.07%  "
Set-TextValue "D8" "3.109.79"
Set-TextValue "E8" "  +0.61%  "
Set-TextValue "E9" "  -0.80%  "
Set-TextValue "E10" "  -2.95%  "
Set-TextValue "D11" "0.152"
Set-TextValue "E11" "  -1.57%  "
Set-TextValue "E12" "  -0.01%  "
Set-TextValue "E13" "  -2.15%  "
Set-TextValue "D14" "37.15"
Set-TextValue "E14" "  +1.18%  "
Set-TextValue "D16" "3.632.47"
Set-TextValue "E16" "  +0.63%  "
Set-TextValue "D17" "66.922.18"
Set-TextValue "E17" "  -0.05%  "
Set-TextValue "D18" "7.16"
Set-TextValue "E18" "  -1.11%  "
Set-TextValue "D19" "3.114.01"
Set-TextValue "E19" "  +0.63%  "
Set-TextValue "D20" "16.30"
Set-TextValue "E20" "  +0.08%  "
Set-TextValue "D21" "476.25"
Set-TextValue "E21" "  +1.49%  "
Set-TextValue "D22" "0.713"
Set-TextValue "E22" "  -0.20%  "
Set-TextValue "D23" "7.92"
Set-TextValue "E23" "  +4.65%  "
Set-TextValue "D24" "13.39"
Set-TextValue "E24" "  +4.01%  "
Set-TextValue "D25" "84.00"
Set-TextValue "E25" "  +0.87%  "
Set-TextValue "D26" "2.28"
Set-TextValue "E26" "  -3.50%  "
Set-TextValue "D27" "10.12"
Set-TextValue "E27" "  -0.53%  "
Set-TextValue "D28" "0.999"
Set-TextValue "E28" "  -0.03%  "
Set-TextValue "D29" "7.87"
Set-TextValue "E29" "  -2.48%  "
Set-TextValue "E30" "  -1.94%  "
Set-TextValue "E31" "  -0.22%  "
Set-TextValue "D32" "28.50"
Set-TextValue "E32" "  +0.90%  "
Set-TextValue "E33" "  -0.21%  "
Set-TextValue "D34" "0.0₃0938"
Set-TextValue "E34" "  -8.80%  "
Set-TextValue "D36" "5.85"
Set-TextValue "E36" "  -0.85%  "
Set-TextValue "D37" "0.971"
Set-TextValue "E37" "  -3.16%  "
Set-TextValue "D38" "46.93"
Set-TextValue "E38" "  +0.38%  "
Set-TextValue "E39" "  -1.24%  "
Set-TextValue "D40" "50.05"
Set-TextValue "E40" "  -0.51%  "
Set-TextValue "D41" "0.310"
Set-TextValue "E41" "  -1.93%  "
Set-TextValue "E42" "  -0.97%  "
Set-TextValue "D43" "8.71"
Set-TextValue "E43" "  -0.04%  "
Set-TextValue "D44" "2.839.65"
Set-TextValue "E44" "  +2.22%  "
Set-TextValue "D45" "386.30"
Set-TextValue "E45" "  -0.13%  "
Set-TextValue "D46" "0.0357"
Set-TextValue "E46" "  -1.59%  "
Set-TextValue "D47" "2.58"
Set-TextValue "E47" "  -9.46%  "
Set-TextValue "E48" "  +0.93%  "
Set-TextValue "E49" "  +0.02%  "
Set-TextValue "D50" "24.86"
Set-TextValue "E50" "  +0.05%  "
Set-TextValue "D51" "2.19"
Set-TextValue "E51" "  -2.37%  "
